# Pre-tx phase self report outcome measures added
# Fills column C ("Pre Experimental Phase") with per-symptom scores for
# rows 2-23, and adds the matching cluster-total formulas in column C
# for rows 24-30 (mirroring the existing column B formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raw per-symptom scores (rows 2-23)
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 5
$ws.Range("C8").Value = 3
$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 1
$ws.Range("C11").Value = 3
$ws.Range("C12").Value = 5
$ws.Range("C13").Value = 3
$ws.Range("C14").Value = 4
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 3
$ws.Range("C19").Value = 3
$ws.Range("C20").Value = 6
$ws.Range("C21").Value = 4
$ws.Range("C22").Value = 6
$ws.Range("C23").Value = 2

# Cluster-total formulas (rows 24-30), mirroring column B's formulas
$ws.Range("C24").Formula = "=SUM(C2:C23)"
$ws.Range("C25").Formula = "=SUM(C2 + C12 + C13)"
$ws.Range("C26").Formula = "=SUM(C19:C22)"
$ws.Range("C27").Formula = "=SUM(C14:C18)"
$ws.Range("C28").Formula = "=C23"
$ws.Range("C29").Formula = "=SUM(C3:C6)"
$ws.Range("C30").Formula = "=SUM(C7:C11)"

# Update the active selection to match the post-edit state
$ws.Range("H27").Select()
